$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D): force text to avoid Excel auto-numeric conversion ---
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '60.119.46'
Set-TextValue $ws.Range('D3') '2.350.67'
Set-TextValue $ws.Range('D4') '0.999'
Set-TextValue $ws.Range('D5') '547.52'
Set-TextValue $ws.Range('D6') '134.02'
Set-TextValue $ws.Range('D7') '0.999'
Set-TextValue $ws.Range('D8') '0.591'
Set-TextValue $ws.Range('D9') '2.348.45'
Set-TextValue $ws.Range('D14') '24.16'
Set-TextValue $ws.Range('D15') '2.765.82'
Set-TextValue $ws.Range('D16') '59.958.99'
Set-TextValue $ws.Range('D18') '2.358.59'
Set-TextValue $ws.Range('D19') '10.74'
Set-TextValue $ws.Range('D21') '317.50'
Set-TextValue $ws.Range('D22') '6.74'
Set-TextValue $ws.Range('D24') '63.12'
Set-TextValue $ws.Range('D26') '0.999'
Set-TextValue $ws.Range('D27') '8.07'
Set-TextValue $ws.Range('D29') '1.75'
Set-TextValue $ws.Range('D31') '171.19'
Set-TextValue $ws.Range('D35') '0.390'
Set-TextValue $ws.Range('D36') '18.02'
Set-TextValue $ws.Range('D40') '317.26'
Set-TextValue $ws.Range('D42') '38.27'
Set-TextValue $ws.Range('D43') '144.91'
Set-TextValue $ws.Range('D45') '0.0960'
Set-TextValue $ws.Range('D47') '0.565'
Set-TextValue $ws.Range('D48') '18.76'
Set-TextValue $ws.Range('D49') '0.0213'
Set-TextValue $ws.Range('D50') '11.02'
Set-TextValue $ws.Range('D51') '1.55'

# Special subscript-three character (U+2083) value for D30 (PEPE ultra-low price)
$sub3 = [char]0x2083
$d30Value = "0.0{0}0757" -f $sub3
Set-TextValue $ws.Range("D30") $d30Value

# --- Volume(1h) column (E): padded percentage text, safe to assign directly ---
$ws.Range('E2').Value = '  +4.00%  '
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('E5').Value = '  +0.92%  '
$ws.Range('E6').Value = '  +2.52%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +4.04%  '
$ws.Range('E9').Value = '  +3.01%  '
$ws.Range('E10').Value = '  +1.39%  '
$ws.Range('E11').Value = '  +2.45%  '
$ws.Range('E12').Value = '  +1.15%  '
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('E14').Value = '  +3.31%  '
$ws.Range('E15').Value = '  +2.96%  '
$ws.Range('E16').Value = '  +3.72%  '
$ws.Range('E17').Value = '  +2.05%  '
$ws.Range('E18').Value = '  +6.31%  '
$ws.Range('E19').Value = '  +2.00%  '
$ws.Range('E20').Value = '  -0.28%  '
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('E22').Value = '  +5.59%  '
$ws.Range('E23').Value = '  +0.44%  '
$ws.Range('E24').Value = '  +0.90%  '
$ws.Range('E25').Value = '  +5.25%  '
$ws.Range('E26').Value = '  -0.56%  '
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('E29').Value = '  +1.29%  '
$ws.Range('E30').Value = '  +5.80%  '
$ws.Range('E31').Value = '  +0.79%  '
$ws.Range('E32').Value = '  +6.94%  '
$ws.Range('E33').Value = '  +4.14%  '
$ws.Range('E34').Value = '  +17.70%  '
$ws.Range('E35').Value = '  +3.05%  '
$ws.Range('E36').Value = '  +2.24%  '
$ws.Range('E37').Value = '  -0.02%  '
$ws.Range('E38').Value = '  -0.04%  '
$ws.Range('E39').Value = '  +6.17%  '
$ws.Range('E40').Value = '  +10.77%  '
$ws.Range('E41').Value = '  +4.79%  '
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('E43').Value = '  +4.17%  '
$ws.Range('E44').Value = '  +2.82%  '
$ws.Range('E45').Value = '  +1.38%  '
$ws.Range('E46').Value = '  +0.58%  '
$ws.Range('E47').Value = '  +2.79%  '
$ws.Range('E48').Value = '  +3.33%  '
$ws.Range('E49').Value = '  +0.90%  '
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('E51').Value = '  +4.44%  '

# --- Coin / Link swap for rows 30 and 31 (Monero <-> PEPE) ---
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('B31').Value = 'Monero'
$ws.Range('C31').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
